# Auto-generated update of Leve price/profit figures across all sheets
# (matches the diff of Sheets/Atomos_Profits.xlsx)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 6250544.5
$ws.Range("I9").Value = 20000172
$ws.Range("J9").Value = 714
$ws.Range("K9").Value = 20000172
$ws.Range("L9").Value = 714
$ws.Range("M9").Value = -20000003
$ws.Range("N9").Value = -1052
# Row 28
$ws.Range("H28").Value = 196
$ws.Range("I28").Value = 209.9
$ws.Range("J28").Value = 168.2
$ws.Range("K28").Value = 209.9
$ws.Range("L28").Value = 168.2
$ws.Range("M28").Value = 275.1
$ws.Range("N28").Value = -1138.2
# Row 92
$ws.Range("H92").Value = 408.9524
$ws.Range("I92").Value = 428.41177
$ws.Range("J92").Value = 326.25
$ws.Range("K92").Value = 428.41177
$ws.Range("L92").Value = 326.25
$ws.Range("M92").Value = 819.5882300000001
$ws.Range("N92").Value = -2822.25
# Row 96
$ws.Range("H96").Value = 736.1177
$ws.Range("I96").Value = 497.3846
$ws.Range("J96").Value = 1512
$ws.Range("K96").Value = 1492.1538
$ws.Range("L96").Value = 4536
$ws.Range("M96").Value = -119.1538
$ws.Range("N96").Value = -7282
# Row 99
$ws.Range("H99").Value = 891.5833
$ws.Range("J99").Value = 1602.4
$ws.Range("L99").Value = 4807.200000000001
$ws.Range("N99").Value = -7803.200000000001
# Row 100
$ws.Range("H100").Value = 2848.258
$ws.Range("I100").Value = 2548.1667
$ws.Range("J100").Value = 3877.1428
$ws.Range("K100").Value = 2548.1667
$ws.Range("L100").Value = 3877.1428
$ws.Range("M100").Value = -2007.1667
$ws.Range("N100").Value = -4959.1428
# Row 101
$ws.Range("H101").Value = 765.0909
$ws.Range("J101").Value = 1691.6666
$ws.Range("L101").Value = 5074.9998
$ws.Range("N101").Value = -8318.9998
# Row 104
$ws.Range("H104").Value = 634.6
$ws.Range("I104").Value = 543.25
$ws.Range("K104").Value = 1629.75
$ws.Range("M104").Value = 117.25
# Row 106
$ws.Range("H106").Value = 875.2727
$ws.Range("I106").Value = 813.55554
$ws.Range("K106").Value = 813.55554
$ws.Range("M106").Value = -182.55554
# Row 112
$ws.Range("H112").Value = 6758083.5
$ws.Range("J112").Value = 7813609
$ws.Range("L112").Value = 23440827
$ws.Range("N112").Value = -23443043

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1775.1724
$ws.Range("I45").Value = 1086.3636
$ws.Range("J45").Value = 3940
$ws.Range("K45").Value = 1086.3636
$ws.Range("L45").Value = 3940
$ws.Range("M45").Value = -709.3635999999999
$ws.Range("N45").Value = -4694
# Row 52
$ws.Range("H52").Value = 39779.375
$ws.Range("J52").Value = 39779.375
$ws.Range("L52").Value = 39779.375
$ws.Range("N52").Value = -40415.375
# Row 61
$ws.Range("H61").Value = 3808.087
$ws.Range("I61").Value = 1633.8334
$ws.Range("J61").Value = 6180
$ws.Range("K61").Value = 1633.8334
$ws.Range("L61").Value = 6180
$ws.Range("M61").Value = -1421.8334
$ws.Range("N61").Value = -6604
# Row 97
$ws.Range("H97").Value = 418
$ws.Range("I97").Value = 446
$ws.Range("J97").Value = 339.6
$ws.Range("K97").Value = 446
$ws.Range("L97").Value = 339.6
$ws.Range("M97").Value = 50
$ws.Range("N97").Value = -1331.6
# Row 102
$ws.Range("H102").Value = 3183.158
$ws.Range("J102").Value = 4995
$ws.Range("L102").Value = 4995
$ws.Range("N102").Value = -8239
# Row 122
$ws.Range("H122").Value = 3187.5
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -17500
# Row 136
$ws.Range("H136").Value = 3808.087
$ws.Range("I136").Value = 1633.8334
$ws.Range("J136").Value = 6180
$ws.Range("K136").Value = 4901.5002
$ws.Range("L136").Value = 18540
$ws.Range("M136").Value = -2351.5002
$ws.Range("N136").Value = -23640

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1004.7083
$ws.Range("J80").Value = 939.2143
$ws.Range("L80").Value = 939.2143
$ws.Range("N80").Value = -2935.2143
# Row 83
$ws.Range("H83").Value = 1004.7083
$ws.Range("J83").Value = 939.2143
$ws.Range("L83").Value = 4696.0715
$ws.Range("N83").Value = -14680.0715
# Row 94
$ws.Range("H94").Value = 295.46155
$ws.Range("I94").Value = 261.9091
$ws.Range("J94").Value = 480
$ws.Range("K94").Value = 261.9091
$ws.Range("L94").Value = 480
$ws.Range("M94").Value = 189.0909
$ws.Range("N94").Value = -1382
# Row 99
$ws.Range("H99").Value = 2772.0667
$ws.Range("I99").Value = 1506.3636
$ws.Range("J99").Value = 6252.75
$ws.Range("K99").Value = 1506.3636
$ws.Range("L99").Value = 6252.75
$ws.Range("M99").Value = -8.363599999999906
$ws.Range("N99").Value = -9248.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2567153.5
$ws.Range("I31").Value = 3227544.2
$ws.Range("J31").Value = 8139.25
$ws.Range("K31").Value = 3227544.2
$ws.Range("L31").Value = 8139.25
$ws.Range("M31").Value = -3227249.2
$ws.Range("N31").Value = -8729.25
# Row 34
$ws.Range("H34").Value = 2567153.5
$ws.Range("I34").Value = 3227544.2
$ws.Range("J34").Value = 8139.25
$ws.Range("K34").Value = 3227544.2
$ws.Range("L34").Value = 8139.25
$ws.Range("M34").Value = -3227342.2
$ws.Range("N34").Value = -8543.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1292
$ws.Range("J5").Value = 2189.4167
$ws.Range("L5").Value = 6568.250100000001
$ws.Range("N5").Value = -6792.250100000001
# Row 16
$ws.Range("H16").Value = 3830.4
$ws.Range("I16").Value = 450
$ws.Range("J16").Value = 4675.5
$ws.Range("K16").Value = 1350
$ws.Range("L16").Value = 14026.5
$ws.Range("M16").Value = -1177
$ws.Range("N16").Value = -14372.5
# Row 135
$ws.Range("H135").Value = 1292
$ws.Range("J135").Value = 2189.4167
$ws.Range("L135").Value = 19704.7503
$ws.Range("N135").Value = -24774.7503

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 50.333332
$ws.Range("I2").Value = 54.25
$ws.Range("J2").Value = 42.5
$ws.Range("K2").Value = 54.25
$ws.Range("L2").Value = 42.5
$ws.Range("M2").Value = 58.75
$ws.Range("N2").Value = -268.5
# Row 97
$ws.Range("H97").Value = 1342.05
$ws.Range("I97").Value = 1030
$ws.Range("J97").Value = 1810.125
$ws.Range("K97").Value = 1030
$ws.Range("L97").Value = 1810.125
$ws.Range("M97").Value = -534
$ws.Range("N97").Value = -2802.125

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 5292.3076
$ws.Range("I16").Value = 2760
$ws.Range("J16").Value = 6875
$ws.Range("K16").Value = 2760
$ws.Range("L16").Value = 6875
$ws.Range("M16").Value = -2590
$ws.Range("N16").Value = -7215
# Row 61
$ws.Range("H61").Value = 90912424
$ws.Range("I61").Value = 166667620
$ws.Range("J61").Value = 6199.8
$ws.Range("K61").Value = 166667620
$ws.Range("L61").Value = 6199.8
$ws.Range("M61").Value = -166667418
$ws.Range("N61").Value = -6603.8
# Row 93
$ws.Range("H93").Value = 1807.7812
$ws.Range("I93").Value = 1306.591
$ws.Range("J93").Value = 2910.4
$ws.Range("K93").Value = 1306.591
$ws.Range("L93").Value = 2910.4
$ws.Range("M93").Value = -58.59099999999989
$ws.Range("N93").Value = -5406.4
# Row 113
$ws.Range("H113").Value = 90912424
$ws.Range("I113").Value = 166667620
$ws.Range("J113").Value = 6199.8
$ws.Range("K113").Value = 166667620
$ws.Range("L113").Value = 6199.8
$ws.Range("M113").Value = -166665450
$ws.Range("N113").Value = -10539.8
# Row 135
$ws.Range("H135").Value = 29470
$ws.Range("J135").Value = 29470
$ws.Range("L135").Value = 29470
$ws.Range("N135").Value = -39610

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 800.6667
$ws.Range("I100").Value = 760.8
$ws.Range("K100").Value = 1521.6
$ws.Range("M100").Value = -980.5999999999999
# Row 138
$ws.Range("H138").Value = 29462
$ws.Range("J138").Value = 29462
$ws.Range("L138").Value = 29462
$ws.Range("N138").Value = -39742

Write-Output "Updated 203 cells across 8 sheets"
